# Adding MPA test automation upload file
# Update the "Data" sheet of the MPA mass-change upload template with the
# new test-automation identifiers (main asset numbers, asset sub-numbers,
# cost center and profit center).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Helper ranges for the numeric-looking values that must stay TEXT (as in
# the original template) rather than be auto-converted to numbers by Excel.
$textCells = @("C6", "C7", "C9", "C11", "D8", "D10", "N7", "N8", "N10")

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# *Main Asset Number (ANLN1)
$ws.Range("C6").Value  = "20000080"
$ws.Range("C7").Value  = "60000369"
$ws.Range("C9").Value  = "60000370"
$ws.Range("C11").Value = "60000371"

# *Asset Subnumber (ANLN2)
$ws.Range("D8").Value  = "268"
$ws.Range("D10").Value = "269"

# Cost Center (KOSTL)
$ws.Range("N7").Value  = "17101904"
$ws.Range("N8").Value  = "17101904"
$ws.Range("N10").Value = "17101904"

# Profit Center (PRCTR) - purely alphanumeric, stays text automatically.
$ws.Range("P7").Value  = "YB103"
$ws.Range("P8").Value  = "YB103"
$ws.Range("P10").Value = "YB103"

# Restore the cells' number format so no residual "Text" formatting is left
# behind on cells that previously used the sheet's default (General) style.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
